$d = $word.ActiveDocument

# The document ends with an empty paragraph. Select it and type the first
# line of text into it, then add a new paragraph with the second line.
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Text = "Hello, this is Nishan Dangal."

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.InsertAfter("I am Studying Managing Software Development unit in this term 4.")
